$wb = $excel.ActiveWorkbook

# Rows (in each localization sheet) that correspond to files whose status is
# "Ready for handoff" and that are getting their handoff info (re)generated.
$rows = @(7, 8, 9, 11, 13, 14)

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-09-05 08:28:42"
}

# --- zh-cn sheet: Priority (E) and Latest Handoff Datetime (H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-09-05 08:28:36"
}

# --- de-de sheet: Priority (E) and Latest Handoff Datetime (H) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-09-05 08:28:42"
}
